# The BOM "R17" (a 2k pulldown resistor for the tach output) was added to the
# CPL placement sheet. In the source data it is physically located right
# before "W1" in the component list, so a whole new row is inserted at row 52
# (pushing every row from the old 52 down to 53, etc.) and populated with the
# new part's placement data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row 52, shifting rows 52:123 down to 53:124.
$ws.Rows.Item(52).Insert()

# Populate the new row with the new component's placement info:
# Designator | Mid X | Mid Y | Layer | Rotation
$ws.Cells.Item(52, 1).Value2 = "R17"
$ws.Cells.Item(52, 2).Value2 = 118.2
$ws.Cells.Item(52, 3).Value2 = -87.5
$ws.Cells.Item(52, 4).Value2 = "top"
$ws.Cells.Item(52, 5).Value2 = 0

# Reflect where the user was last working when the file was saved.
[void]$ws.Range("D52").Select()
$excel.ActiveWindow.ScrollRow = 23
